# Insert one new data row at row 122 (pushing the previous row 122 and all
# following rows down by one, through the former row 190 which becomes 191),
# then populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(122).Insert()

$ws.Cells.Item(122, 1).Value = 10
$ws.Cells.Item(122, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(122, 3).Value = "La Araucanía"
$ws.Cells.Item(122, 4).Value = 44488
$ws.Cells.Item(122, 5).Value = 9
$ws.Cells.Item(122, 6).Value = 100112017
$ws.Cells.Item(122, 7).Value = "Apio"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 235
$ws.Cells.Item(122, 11).Value = 8000
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = 8468
$ws.Cells.Item(122, 14).Value = "`$/docena de matas"
$ws.Cells.Item(122, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(122, 16).Value = 1411
$ws.Cells.Item(122, 17).Value = 6
$ws.Cells.Item(122, 18).Value = "Hortaliza"
